$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "Unordered-DenStream"
$ws.Range("D1").Select()
